$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-looking numeric strings (e.g. "0.999",
# "7.77") as text. Force text format on each Price cell we touch so COM
# does not silently convert the literal into a Number cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.616.95'
$ws.Range("E2").Value = '  +1.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.452.01'
$ws.Range("E3").Value = '  +1.81%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.30'
$ws.Range("E5").Value = '  +1.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.21'
$ws.Range("E6").Value = '  +8.81%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.454.13'
$ws.Range("E7").Value = '  +1.87%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.475'
$ws.Range("E9").Value = '  +0.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.83'
$ws.Range("E10").Value = '  +2.42%  '

$ws.Range("E11").Value = '  +1.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.392'
$ws.Range("E12").Value = '  +1.64%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.041.00'
$ws.Range("E13").Value = '  +1.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.98'
$ws.Range("E14").Value = '  +6.93%  '

$ws.Range("E15").Value = '  -0.35%  '

$ws.Range("E16").Value = '  +1.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.446.69'
$ws.Range("E17").Value = '  +1.66%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.714.49'
$ws.Range("E18").Value = '  +0.95%  '

$ws.Range("E19").Value = '  +8.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.41'
$ws.Range("E20").Value = '  +2.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.52'
$ws.Range("E21").Value = '  +0.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.75'
$ws.Range("E22").Value = '  +3.54%  '

$ws.Range("E23").Value = '  +2.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.585.80'
$ws.Range("E24").Value = '  +1.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.94'
$ws.Range("E25").Value = '  +2.49%  '

$ws.Range("E26").Value = '  -0.19%  '

$ws.Range("E27").Value = '  +0.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000125'
$ws.Range("E28").Value = '  -1.20%  '

$ws.Range("E29").Value = '  +3.79%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.77'
$ws.Range("E30").Value = '  +3.21%  '

$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.55'
$ws.Range("E31").Value = '  -12.60%  '

$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  +0.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.27'
$ws.Range("E33").Value = '  +1.26%  '

$ws.Range("E34").Value = '  +0.91%  '

$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '24.04'
$ws.Range("E36").Value = '  +1.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.26'
$ws.Range("E37").Value = '  +0.78%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.07'
$ws.Range("E38").Value = '  +2.81%  '

$ws.Range("E39").Value = '  +1.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '166.96'
$ws.Range("E40").Value = '  +1.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0790'
$ws.Range("E41").Value = '  +3.87%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.15'
$ws.Range("E42").Value = '  +12.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.795'
$ws.Range("E43").Value = '  +1.91%  '

$ws.Range("E44").Value = '  +2.24%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.33'
$ws.Range("E46").Value = '  +1.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.72'
$ws.Range("E47").Value = '  +0.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.603.58'
$ws.Range("E48").Value = '  +5.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.17'
$ws.Range("E49").Value = '  -2.77%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.00'
$ws.Range("E50").Value = '  +2.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.16'
$ws.Range("E51").Value = '  -0.02%  '
